# Fix price/volume/ticker data that had been incorrectly filled in from
# other companies' rows (extra files bug). Rows 2-14 and 16-25 get their
# open/close/high/low price, shares outstanding and fixed_ticker columns
# corrected to reference CGNT (Cognyte Software) instead of unrelated
# tickers (BABA, MRVL, UBER, MU, NFLX, ...). Row 15 already held the
# correct CGNT values and is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = 27.98999977111816
$ws.Range("E2").Value = 28.8700008392334
$ws.Range("F2").Value = 32.45999908447266
$ws.Range("G2").Value = 26.77000045776367
$ws.Range("H2").Value = 72969110
$ws.Range("I2").Value = "CGNT"
$ws.Range("D3").Value = 27.98999977111816
$ws.Range("E3").Value = 28.8700008392334
$ws.Range("F3").Value = 32.45999908447266
$ws.Range("G3").Value = 26.77000045776367
$ws.Range("H3").Value = 72969110
$ws.Range("I3").Value = "CGNT"
$ws.Range("D4").Value = 27.98999977111816
$ws.Range("E4").Value = 28.8700008392334
$ws.Range("F4").Value = 32.45999908447266
$ws.Range("G4").Value = 26.77000045776367
$ws.Range("H4").Value = 72969110
$ws.Range("I4").Value = "CGNT"
$ws.Range("D5").Value = 27.98999977111816
$ws.Range("E5").Value = 28.8700008392334
$ws.Range("F5").Value = 32.45999908447266
$ws.Range("G5").Value = 26.77000045776367
$ws.Range("H5").Value = 72969110
$ws.Range("I5").Value = "CGNT"
$ws.Range("D6").Value = 27.98999977111816
$ws.Range("E6").Value = 28.8700008392334
$ws.Range("F6").Value = 32.45999908447266
$ws.Range("G6").Value = 26.77000045776367
$ws.Range("H6").Value = 72969110
$ws.Range("I6").Value = "CGNT"
$ws.Range("D7").Value = 27.98999977111816
$ws.Range("E7").Value = 28.8700008392334
$ws.Range("F7").Value = 32.45999908447266
$ws.Range("G7").Value = 26.77000045776367
$ws.Range("H7").Value = 72969110
$ws.Range("I7").Value = "CGNT"
$ws.Range("D8").Value = 26.34000015258789
$ws.Range("E8").Value = 25.73999977111816
$ws.Range("F8").Value = 26.54000091552734
$ws.Range("G8").Value = 22.61000061035156
$ws.Range("H8").Value = 72969110
$ws.Range("I8").Value = "CGNT"
$ws.Range("D9").Value = 26.03000068664551
$ws.Range("E9").Value = 27.47999954223633
$ws.Range("F9").Value = 28.27000045776367
$ws.Range("G9").Value = 26.03000068664551
$ws.Range("H9").Value = 72969110
$ws.Range("I9").Value = "CGNT"
$ws.Range("D10").Value = 19.92000007629395
$ws.Range("E10").Value = 20.3700008392334
$ws.Range("F10").Value = 24.23999977111816
$ws.Range("G10").Value = 19.70000076293945
$ws.Range("H10").Value = 72969110
$ws.Range("I10").Value = "CGNT"
$ws.Range("D11").Value = 10.92000007629394
$ws.Range("E11").Value = 11.15999984741211
$ws.Range("F11").Value = 11.38799953460693
$ws.Range("G11").Value = 9.300000190734863
$ws.Range("H11").Value = 72969110
$ws.Range("I11").Value = "CGNT"
$ws.Range("D12").Value = 6.760000228881836
$ws.Range("E12").Value = 7.079999923706055
$ws.Range("F12").Value = 7.590000152587891
$ws.Range("G12").Value = 5.960000038146973
$ws.Range("H12").Value = 72969110
$ws.Range("I12").Value = "CGNT"
$ws.Range("D13").Value = 4.5
$ws.Range("E13").Value = 5.21999979019165
$ws.Range("F13").Value = 5.340000152587891
$ws.Range("G13").Value = 4.440000057220459
$ws.Range("H13").Value = 72969110
$ws.Range("I13").Value = "CGNT"
$ws.Range("D14").Value = 2.759999990463257
$ws.Range("E14").Value = 2.950000047683716
$ws.Range("F14").Value = 3.700000047683716
$ws.Range("G14").Value = 2.575000047683716
$ws.Range("H14").Value = 72969110
$ws.Range("I14").Value = "CGNT"
$ws.Range("D16").Value = 4.159999847412109
$ws.Range("E16").Value = 4.96999979019165
$ws.Range("F16").Value = 4.96999979019165
$ws.Range("G16").Value = 4.025000095367432
$ws.Range("H16").Value = 72969110
$ws.Range("I16").Value = "CGNT"
$ws.Range("D17").Value = 5.440000057220459
$ws.Range("E17").Value = 4.789999961853027
$ws.Range("F17").Value = 5.789999961853027
$ws.Range("G17").Value = 4.579999923706055
$ws.Range("H17").Value = 72969110
$ws.Range("I17").Value = "CGNT"
$ws.Range("D18").Value = 4.210000038146973
$ws.Range("E18").Value = 5.130000114440918
$ws.Range("F18").Value = 5.429999828338623
$ws.Range("G18").Value = 4.010000228881836
$ws.Range("H18").Value = 72969110
$ws.Range("I18").Value = "CGNT"
$ws.Range("D19").Value = 7.150000095367432
$ws.Range("E19").Value = 7.420000076293945
$ws.Range("F19").Value = 7.78000020980835
$ws.Range("G19").Value = 6.829999923706055
$ws.Range("H19").Value = 72969110
$ws.Range("I19").Value = "CGNT"
$ws.Range("D20").Value = 6.809999942779541
$ws.Range("E20").Value = 7.599999904632568
$ws.Range("F20").Value = 7.869999885559082
$ws.Range("G20").Value = 6.769999980926514
$ws.Range("H20").Value = 72969110
$ws.Range("I20").Value = "CGNT"
$ws.Range("D21").Value = 7.659999847412109
$ws.Range("E21").Value = 7.659999847412109
$ws.Range("F21").Value = 7.690000057220459
$ws.Range("G21").Value = 6.739999771118164
$ws.Range("H21").Value = 72969110
$ws.Range("I21").Value = "CGNT"
$ws.Range("D22").Value = 6.5
$ws.Range("E22").Value = 8.119999885559082
$ws.Range("F22").Value = 8.239999771118164
$ws.Range("G22").Value = 6.380000114440918
$ws.Range("H22").Value = 72969110
$ws.Range("I22").Value = "CGNT"
$ws.Range("D23").Value = 9.029999732971191
$ws.Range("E23").Value = 8.930000305175781
$ws.Range("F23").Value = 11.11999988555908
$ws.Range("G23").Value = 8.630000114440918
$ws.Range("H23").Value = 72969110
$ws.Range("I23").Value = "CGNT"
$ws.Range("D24").Value = 9.90999984741211
$ws.Range("E24").Value = 10.90499973297119
$ws.Range("F24").Value = 11.09000015258789
$ws.Range("G24").Value = 9.579999923706056
$ws.Range("H24").Value = 72969110
$ws.Range("I24").Value = "CGNT"
$ws.Range("D25").Value = 9.100000381469728
$ws.Range("E25").Value = 8.930000305175781
$ws.Range("F25").Value = 9.350000381469728
$ws.Range("G25").Value = 8.34000015258789
$ws.Range("H25").Value = 72969110
$ws.Range("I25").Value = "CGNT"
